$d = $word.ActiveDocument
$w = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# ---------------------------------------------------------------
# 1. "Features" section heading: remove the blank paragraph before
#    "Features" and promote "Features" itself to Heading1.
# ---------------------------------------------------------------
$d.Paragraphs.Item(8).Range.Delete()          # was blank spacer para
$d.Paragraphs.Item(8).Style = "Heading1"       # "Features" -> Heading1

# ---------------------------------------------------------------
# 2. New intro paragraph right after the "Features" heading.
# ---------------------------------------------------------------
$d.Paragraphs.Item(8).Range.InsertParagraphAfter()
$d.Paragraphs.Item(9).Range.InsertXML('<w:p xmlns:w="' + $w + '"><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>The Progress Tree Control has the following features:</w:t></w:r></w:p>')

# ---------------------------------------------------------------
# 3. Rewrite the 3 top-level bullet texts under "Features".
# ---------------------------------------------------------------
$d.Paragraphs.Item(10).Range.Find.Execute("Allows selecting a Progress element or the entire tree (by selecting the top element)", $true, $false, $false, $false, $false, $true, 1, $false, "The individual elements of the tree can be selected. Upon selection an event will be fired indicating which element(s) have been selected. Programmatically selecting an element is also possible and also leads to the firing of the selection event.", 2) | Out-Null
$d.Paragraphs.Item(11).Range.Find.Execute("Only allows a single top-level node", $true, $false, $false, $false, $false, $true, 1, $false, "There can only be one top-level node.", 2) | Out-Null
$d.Paragraphs.Item(12).Range.Find.Execute("Shows:", $true, $false, $false, $false, $false, $true, 1, $false, "Each element of the Progress Tree Control contains:", 2) | Out-Null

# ---------------------------------------------------------------
# 4. "Time estimate" / "Text indicating current action" rewrites.
# ---------------------------------------------------------------
$d.Paragraphs.Item(14).Range.Find.Execute("Time estimate", $true, $false, $false, $false, $false, $true, 1, $false, "An estimate for the remaining time", 2) | Out-Null
$d.Paragraphs.Item(15).Range.Find.Execute("Text indicating current action", $true, $false, $false, $false, $false, $true, 1, $false, "An indication of the current action that is being performed", 2) | Out-Null

# ---------------------------------------------------------------
# 5. Merge "Buttons to start / stop / pause ..." (para 16) and
#    "Expanding section with more details" (para 17) into a single
#    multi-run paragraph; delete the now-redundant paragraph.
# ---------------------------------------------------------------
$d.Paragraphs.Item(17).Range.Delete()
$p16 = $d.Paragraphs.Item(16)
$p16.Range.InsertXML('<w:p xmlns:w="' + $w + '"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>S</w:t></w:r><w:r><w:t>tart</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> stop </w:t></w:r><w:r><w:t>and</w:t></w:r><w:r><w:t xml:space="preserve"> pause</w:t></w:r><w:r><w:t xml:space="preserve"> buttons which action on the </w:t></w:r><w:r><w:t>action which is being monitored</w:t></w:r></w:p>')

# ---------------------------------------------------------------
# 6. "Must provide events for:" -> two runs ending in a colon.
# ---------------------------------------------------------------
$d.Paragraphs.Item(17).Range.InsertXML('<w:p xmlns:w="' + $w + '"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>The Progress Tree Control provides the following events</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p>')

# ---------------------------------------------------------------
# 7. "OnSelect: provides selected item" -> "SelectedItemChanged"
#    (drop the trailing ": provides selected item" run).
# ---------------------------------------------------------------
$d.Paragraphs.Item(18).Range.InsertXML('<w:p xmlns:w="' + $w + '"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>SelectedItemChanged</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

Write-Host "Stage7 count:" $d.Paragraphs.Count

# ---------------------------------------------------------------
# 8. Tail of the document: new "Automatically generated ..." bullet,
#    four blank spacer paragraphs, and the whole new "Design" section.
#    (The very last, already-existing blank paragraph is left as is.)
# ---------------------------------------------------------------
$designXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Automatically generated children cannot be stopped / paused. In order to stop an automatically generated child the parent action must be stopped.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:after="0"/></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Design</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">Develop separate component for </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>TreeNode.Header</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> --&gt; </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Has</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>progressbar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> etc.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Tree actions:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Add new tree element</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Remove tree element</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Select tree element</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>TreeNode</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> actions:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Update progress</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Update current action</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Show details</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Pause action</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t>Stop action (--&gt; leads to removal of the node)</w:t></w:r></w:p>
'@
$tailAnchor = $d.Paragraphs.Item($d.Paragraphs.Count - 1).Range
$tailAnchor.Collapse(1)
$tailAnchor.InsertXML($designXml)

Write-Host "DEBUG after insert, count:" $d.Paragraphs.Count
for ($i = 26; $i -le 34; $i++) {
    Write-Host "DEBUG" $i ":: [" $d.Paragraphs.Item($i).Range.Text "] style=" $d.Paragraphs.Item($i).Style.NameLocal
}

# Remove the stray blank paragraph that used to sit right after "Error"
# (it has now been pushed down, just before the newly-inserted block).
$d.Paragraphs.Item(29).Range.Delete()

Write-Host "Final count:" $d.Paragraphs.Count
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Host $i ":: [" $d.Paragraphs.Item($i).Range.Text "] style=" $d.Paragraphs.Item($i).Style.NameLocal
}
